# Updates the Price (D) and Volume(1h) (E) columns for the cryptos list,
# reflecting a refreshed data pull (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column letter, new text value.
$updates = @(
    @{ Row = 2; Col = 'D'; Value = '25.889.26' }
    @{ Row = 2; Col = 'E'; Value = '  -3.38%  ' }
    @{ Row = 3; Col = 'D'; Value = '1.828.46' }
    @{ Row = 3; Col = 'E'; Value = '  -2.17%  ' }
    @{ Row = 4; Col = 'D'; Value = '0.9980' }
    @{ Row = 4; Col = 'E'; Value = '  -0.28%  ' }
    @{ Row = 5; Col = 'D'; Value = '278.11' }
    @{ Row = 5; Col = 'E'; Value = '  -7.31%  ' }
    @{ Row = 6; Col = 'D'; Value = '0.9981' }
    @{ Row = 6; Col = 'E'; Value = '  -0.27%  ' }
    @{ Row = 7; Col = 'D'; Value = '0.5104' }
    @{ Row = 7; Col = 'E'; Value = '  -4.33%  ' }
    @{ Row = 8; Col = 'E'; Value = '  -7.00%  ' }
    @{ Row = 9; Col = 'D'; Value = '44.67' }
    @{ Row = 9; Col = 'E'; Value = '  -1.56%  ' }
    @{ Row = 10; Col = 'D'; Value = '0.06792' }
    @{ Row = 10; Col = 'E'; Value = '  -4.87%  ' }
    @{ Row = 11; Col = 'D'; Value = '19.84' }
    @{ Row = 11; Col = 'E'; Value = '  -7.82%  ' }
    @{ Row = 12; Col = 'D'; Value = '0.8079' }
    @{ Row = 12; Col = 'E'; Value = '  -8.88%  ' }
    @{ Row = 13; Col = 'D'; Value = '0.07811' }
    @{ Row = 13; Col = 'E'; Value = '  -4.22%  ' }
    @{ Row = 14; Col = 'D'; Value = '1.827.95' }
    @{ Row = 14; Col = 'E'; Value = '  -2.35%  ' }
    @{ Row = 15; Col = 'D'; Value = '5.071' }
    @{ Row = 15; Col = 'E'; Value = '  -4.10%  ' }
    @{ Row = 16; Col = 'D'; Value = '87.97' }
    @{ Row = 16; Col = 'E'; Value = '  -4.74%  ' }
    @{ Row = 17; Col = 'D'; Value = '0.9980' }
    @{ Row = 17; Col = 'E'; Value = '  -0.26%  ' }
    @{ Row = 18; Col = 'D'; Value = '14.14' }
    @{ Row = 18; Col = 'E'; Value = '  -4.70%  ' }
    @{ Row = 19; Col = 'D'; Value = '0.000008059' }
    @{ Row = 19; Col = 'E'; Value = '  -4.96%  ' }
    @{ Row = 20; Col = 'D'; Value = '0.9976' }
    @{ Row = 20; Col = 'E'; Value = '  -0.30%  ' }
    @{ Row = 21; Col = 'D'; Value = '25.930.97' }
    @{ Row = 21; Col = 'E'; Value = '  -3.33%  ' }
    @{ Row = 22; Col = 'D'; Value = '4.768' }
    @{ Row = 22; Col = 'E'; Value = '  -4.02%  ' }
    @{ Row = 23; Col = 'E'; Value = '  -5.81%  ' }
    @{ Row = 24; Col = 'D'; Value = '6.180' }
    @{ Row = 24; Col = 'E'; Value = '  -3.07%  ' }
    @{ Row = 25; Col = 'D'; Value = '2.356' }
    @{ Row = 25; Col = 'E'; Value = '  +3.05%  ' }
    @{ Row = 26; Col = 'D'; Value = '142.62' }
    @{ Row = 26; Col = 'E'; Value = '  -2.37%  ' }
    @{ Row = 27; Col = 'D'; Value = '1.667' }
    @{ Row = 27; Col = 'E'; Value = '  -4.14%  ' }
    @{ Row = 28; Col = 'D'; Value = '17.19' }
    @{ Row = 29; Col = 'D'; Value = '109.35' }
    @{ Row = 29; Col = 'E'; Value = '  -3.82%  ' }
    @{ Row = 30; Col = 'D'; Value = '4.337' }
    @{ Row = 30; Col = 'E'; Value = '  -7.64%  ' }
    @{ Row = 31; Col = 'D'; Value = '4.291' }
    @{ Row = 31; Col = 'E'; Value = '  -7.19%  ' }
    @{ Row = 32; Col = 'D'; Value = '0.08755' }
    @{ Row = 32; Col = 'E'; Value = '  -3.73%  ' }
    @{ Row = 33; Col = 'D'; Value = '0.04864' }
    @{ Row = 33; Col = 'E'; Value = '  -3.07%  ' }
    @{ Row = 34; Col = 'D'; Value = '1.167' }
    @{ Row = 34; Col = 'E'; Value = '  -0.19%  ' }
    @{ Row = 35; Col = 'D'; Value = '0.7283' }
    @{ Row = 35; Col = 'E'; Value = '  -9.79%  ' }
    @{ Row = 36; Col = 'D'; Value = '2.864' }
    @{ Row = 36; Col = 'E'; Value = '  -2.64%  ' }
    @{ Row = 37; Col = 'D'; Value = '3.188' }
    @{ Row = 37; Col = 'E'; Value = '  +0.01%  ' }
    @{ Row = 38; Col = 'D'; Value = '2.396' }
    @{ Row = 39; Col = 'D'; Value = '0.01848' }
    @{ Row = 39; Col = 'E'; Value = '  -4.86%  ' }
    @{ Row = 40; Col = 'D'; Value = '0.5112' }
    @{ Row = 40; Col = 'E'; Value = '  -16.27%  ' }
    @{ Row = 41; Col = 'D'; Value = '0.9460' }
    @{ Row = 41; Col = 'E'; Value = '  -11.16%  ' }
    @{ Row = 42; Col = 'D'; Value = '117.03' }
    @{ Row = 42; Col = 'E'; Value = '  +0.72%  ' }
    @{ Row = 43; Col = 'D'; Value = '6.217' }
    @{ Row = 43; Col = 'E'; Value = '  -4.07%  ' }
    @{ Row = 44; Col = 'D'; Value = '8.005' }
    @{ Row = 44; Col = 'E'; Value = '  -8.62%  ' }
    @{ Row = 45; Col = 'D'; Value = '0.9970' }
    @{ Row = 45; Col = 'E'; Value = '  -0.35%  ' }
    @{ Row = 46; Col = 'D'; Value = '0.1361' }
    @{ Row = 46; Col = 'E'; Value = '  -8.37%  ' }
    @{ Row = 47; Col = 'D'; Value = '0.4493' }
    @{ Row = 47; Col = 'E'; Value = '  -14.72%  ' }
    @{ Row = 48; Col = 'D'; Value = '9.342' }
    @{ Row = 48; Col = 'E'; Value = '  -6.06%  ' }
    @{ Row = 49; Col = 'E'; Value = '  -2.76%  ' }
    @{ Row = 50; Col = 'D'; Value = '0.05920' }
    @{ Row = 50; Col = 'E'; Value = '  -2.41%  ' }
    @{ Row = 51; Col = 'D'; Value = '1.493' }
    @{ Row = 51; Col = 'E'; Value = '  -9.27%  ' }
)

$colIndex = @{ D = 4; E = 5 }

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $colIndex[$u.Col])
    # Force text so numeric-looking strings (e.g. "278.11") are not
    # coerced into floating-point numbers, and restore the default
    # "Normal" style afterwards so no stray number-format style lingers.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}
